# "Pompe feature P1 & P3"
#
# 1. "map (4)" sheet: flag columns F:H for rows 6-8 (pumps P1 & P3) as
#    active (0 -> 1), then leave the selection on J8.
# 2. "confs" sheet: fill in the computed coefficients for row 3 (same
#    figures/formatting as row 2's H:J columns) and leave the selection on
#    G3 - this also becomes the active tab when the workbook is saved,
#    taking over from "map (2)".

$wb = $excel.ActiveWorkbook

# --- "map (4)": turn on P1 & P3 flags for rows 6-8 ---
$wsMap4 = $wb.Worksheets.Item("map (4)")
$wsMap4.Range("F6:H8").Value = 1
[void]$wsMap4.Range("J8").Select()

# --- "confs": row 3 picks up the same style + values as row 2's H:J ---
$wsConfs = $wb.Worksheets.Item("confs")
[void]$wsConfs.Range("H2:J2").Copy()
[void]$wsConfs.Range("H3:J3").PasteSpecial(-4122)
$wsConfs.Range("H3").Value = -0.000516
$wsConfs.Range("I3").Value = -0.0154
$wsConfs.Range("J3").Value = 4.87

# Last selection made wins for both "which sheet is active" and the
# sheet's own stored selection, so do this last to make "confs" active.
[void]$wsConfs.Range("G3").Select()
